$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price/volume snapshot data. Price cells in column D are
# stored as text in the source sheet (e.g. "69.282.79", "7.00"), so force
# the Text number format before assigning numeric-looking strings to keep
# them from being reinterpreted as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.321.75'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.682.34'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '682.78'
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.25'
$ws.Range("E6").Value = '  -2.65%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -1.13%  '
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.00'
$ws.Range("E10").Value = '  -3.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.437'
$ws.Range("E11").Value = '  -3.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000233'
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.301.53'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.26'
$ws.Range("E14").Value = '  -3.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.681.48'
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.318.73'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("E17").Value = '  +1.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.97'
$ws.Range("E18").Value = '  -2.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.39'
$ws.Range("E19").Value = '  -4.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '470.11'
$ws.Range("E20").Value = '  -1.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.04'
$ws.Range("E21").Value = '  +2.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.650'
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '80.01'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.825.44'
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000121'
$ws.Range("E26").Value = '  -5.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.91'
$ws.Range("E27").Value = '  -5.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.15'
$ws.Range("E28").Value = '  -4.88%  '
$ws.Range("E29").Value = '  -1.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.75'
$ws.Range("E30").Value = '  -4.93%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.57'
$ws.Range("E32").Value = '  -4.76%  '
$ws.Range("E33").Value = '  -5.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.85'
$ws.Range("E34").Value = '  -1.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.660.26'
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("E36").Value = '  -3.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.20'
$ws.Range("E37").Value = '  -4.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.06'
$ws.Range("E38").Value = '  -1.28%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.23'
$ws.Range("E40").Value = '  +3.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0903'
$ws.Range("E41").Value = '  -4.18%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '168.02'
$ws.Range("E43").Value = '  +8.90%  '
$ws.Range("E44").Value = '  -1.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.62'
$ws.Range("E46").Value = '  -4.74%  '
$ws.Range("E47").Value = '  -1.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.11'
$ws.Range("E48").Value = '  +2.79%  '
$ws.Range("E49").Value = '  -5.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.77'
$ws.Range("E50").Value = '  -4.18%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.267'
$ws.Range("E51").Value = '  -2.14%  '
